{"js": "// Update the date line and the 25 two-digit multiplication problems in the\n// table to the new values from the target revision. All old values in this\n// document are unique strings, so a plain text search-and-replace (via\n// Body.search + Range.insertText) is sufficient and keeps each run's\n// existing formatting (font/size) intact.\nconst replacements = [\n  [\"2025-09-09 Tuesday\", \"2025-09-10 Wednesday\"],\n  [\"96\u00d760=\", \"89\u00d735=\"],\n  [\"41\u00d728=\", \"97\u00d723=\"],\n  [\"15\u00d799=\", \"33\u00d751=\"],\n  [\"70\u00d752=\", \"91\u00d717=\"],\n  [\"21\u00d719=\", \"42\u00d745=\"],\n  [\"97\u00d716=\", \"13\u00d749=\"],\n  [\"60\u00d787=\", \"85\u00d721=\"],\n  [\"87\u00d741=\", \"42\u00d712=\"],\n  [\"86\u00d773=\", \"19\u00d738=\"],\n  [\"87\u00d749=\", \"35\u00d717=\"],\n  [\"58\u00d755=\", \"70\u00d788=\"],\n  [\"65\u00d755=\", \"67\u00d743=\"],\n  [\"55\u00d767=\", \"76\u00d788=\"],\n  [\"62\u00d755=\", \"63\u00d752=\"],\n  [\"90\u00d773=\", \"65\u00d783=\"],\n  [\"77\u00d757=\", \"12\u00d786=\"],\n  [\"63\u00d799=\", \"38\u00d733=\"],\n  [\"62\u00d748=\", \"55\u00d716=\"],\n  [\"58\u00d751=\", \"86\u00d783=\"],\n  [\"46\u00d793=\", \"36\u00d711=\"],\n  [\"98\u00d727=\", \"81\u00d729=\"],\n  [\"69\u00d798=\", \"18\u00d798=\"],\n  [\"55\u00d748=\", \"43\u00d752=\"],\n  [\"11\u00d761=\", \"93\u00d724=\"],\n  [\"60\u00d799=\", \"32\u00d764=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 two-digit multiplication problems in the\n# table to the new values from the target revision. All old values in this\n# document are unique strings, so a plain Find/Replace over the whole\n# document body is sufficient and preserves each run's existing formatting\n# (font/size), since Find.Execute only rewrites the matched text.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-09-09 Tuesday\", \"2025-09-10 Wednesday\"),\n    @(\"96\u00d760=\", \"89\u00d735=\"),\n    @(\"41\u00d728=\", \"97\u00d723=\"),\n    @(\"15\u00d799=\", \"33\u00d751=\"),\n    @(\"70\u00d752=\", \"91\u00d717=\"),\n    @(\"21\u00d719=\", \"42\u00d745=\"),\n    @(\"97\u00d716=\", \"13\u00d749=\"),\n    @(\"60\u00d787=\", \"85\u00d721=\"),\n    @(\"87\u00d741=\", \"42\u00d712=\"),\n    @(\"86\u00d773=\", \"19\u00d738=\"),\n    @(\"87\u00d749=\", \"35\u00d717=\"),\n    @(\"58\u00d755=\", \"70\u00d788=\"),\n    @(\"65\u00d755=\", \"67\u00d743=\"),\n    @(\"55\u00d767=\", \"76\u00d788=\"),\n    @(\"62\u00d755=\", \"63\u00d752=\"),\n    @(\"90\u00d773=\", \"65\u00d783=\"),\n    @(\"77\u00d757=\", \"12\u00d786=\"),\n    @(\"63\u00d799=\", \"38\u00d733=\"),\n    @(\"62\u00d748=\", \"55\u00d716=\"),\n    @(\"58\u00d751=\", \"86\u00d783=\"),\n    @(\"46\u00d793=\", \"36\u00d711=\"),\n    @(\"98\u00d727=\", \"81\u00d729=\"),\n    @(\"69\u00d798=\", \"18\u00d798=\"),\n    @(\"55\u00d748=\", \"43\u00d752=\"),\n    @(\"11\u00d761=\", \"93\u00d724=\"),\n    @(\"60\u00d799=\", \"32\u00d764=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
